# Apply Greek translations to GET SMARTCASH.docx
$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $found = $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false)
    if (-not $found) {
        throw "Could not find text: $old"
    }
    $range.Text = $new
}

# 1. GPU NVIDIA miner -> GPU NVIDIA Miner (Κάρτας γραφικών Νvidia)
Replace-Text "GPU NVIDIA miner" "GPU NVIDIA Miner (Κάρτας γραφικών Νvidia)"

# 2. GPU AMD miner -> GPU AMD Miner (Κάρτας γραφικών AMD)
Replace-Text "GPU AMD miner" "GPU AMD Miner (Κάρτας γραφικών AMD)"

# 3. For mining support please join -> Για βοήθεια για το Mining μπείτε στο:
Replace-Text "For mining support please join" "Για βοήθεια για το Mining μπείτε στο:"

# 4. Discord -> " Discord" (add a leading space)
Replace-Text "Discord" " Discord"

# 5. " EXCHANGES" -> "Ανταλλακτήρια"
Replace-Text " EXCHANGES" "Ανταλλακτήρια"

# 6. Turn altcoins into Smartcash instantly -> Μετατρέψει το altcoins σε Smartcash αμέσως
Replace-Text "Turn altcoins into Smartcash instantly" "Μετατρέψει το altcoins σε Smartcash αμέσως"
